$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2021-10-15"

# Update the row label for the October partial-month row
$ws.Range("A12").Value = "October (through 10-15)"

# Row 12 (October, partial month) updates
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 23
$ws.Range("G12").Value = 0.1154
$ws.Range("H12").Value = 6
$ws.Range("I12").Value = 23
$ws.Range("J12").Value = 0.2069
$ws.Range("L12").Value = 36
$ws.Range("M12").Value = 0.0526
$ws.Range("R12").Value = 75
$ws.Range("U12").Value = 95

# Row 13 (Total) updates
$ws.Range("E13").Value = 49
$ws.Range("F13").Value = 406
$ws.Range("G13").Value = 0.1077
$ws.Range("H13").Value = 56
$ws.Range("I13").Value = 600
$ws.Range("J13").Value = 0.0854
$ws.Range("L13").Value = 523
$ws.Range("M13").Value = 0.1075
$ws.Range("R13").Value = 923
$ws.Range("S13").Value = 0.0543
$ws.Range("U13").Value = 1263
$ws.Range("V13").Value = 0.061
